$p = $ppt.ActivePresentation

# Slide 1 title: "First" + " " + "slide" -> single run "First slide"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Delete()
$tr1.InsertAfter("First slide")

# The notes page bound to slide 2 in the package (ppt/notesSlides/notesSlide1.xml,
# linked via slide2.xml.rels) holds the speaker notes text for the deck's first
# slide. Merge its per-word runs into a single run.
$s2 = $p.Slides.Item(2)
$notesTr = $s2.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesTr.Text = "Some notes here: this first slide should use the Blank template"

# Slide 3 title: "Third" + " " + "slide" -> single run "Third slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Delete()
$tr3.InsertAfter("Third slide")
